$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: B, C, D, E, F, I, J, K, L, M, N (column indices 2,3,4,5,6,9,10,11,12,13,14)
$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$rowData = @{
    2 = @(1.02, 1.046220375459431, 1.053886849242428, 1.053598608582881, 1.06376747284704, 1.043905258353602, 1.051275932400249, 1.056631579472896, 1.056344134153944, 1.066485221650932, 1.020904100012277)
    3 = @(1.02, 1.047241206148459, 1.054557949145935, 1.054475845473333, 1.06465510024825, 1.044134300124469, 1.051944391112135, 1.05711617674884, 1.057034283290932, 1.067187741855793, 1.021130675458881)
    4 = @(1.02, 1.047902230417211, 1.054992548558462, 1.055044245255549, 1.065230253682419, 1.044281531998373, 1.052376802550586, 1.057429413212578, 1.057480984194527, 1.067642478163081, 1.021277133313038)
    5 = @(1.02, 1.048180238968567, 1.055175337234724, 1.055283383420295, 1.065472238174725, 1.044343194729999, 1.052558557356044, 1.057561017761279, 1.0576688070383, 1.067833686230574, 1.021338667540455)
    6 = @(1.02, 1.048226924491838, 1.055206033089039, 1.055323546461591, 1.06551287953815, 1.044353534464155, 1.052589072949109, 1.057583110021181, 1.057700345039996, 1.067865793040237, 1.021348997248149)
    7 = @(1.02, 1.047905944734097, 1.054994990665003, 1.055047439914216, 1.06523348634619, 1.044282356857114, 1.052379231289838, 1.057431172034795, 1.057483493776681, 1.067645032950908, 1.021277955680871)
    8 = @(1.02, 1.046565270510025, 1.054113576684417, 1.05389491499745, 1.064067284413697, 1.043982865355286, 1.05150186677634, 1.056795419216353, 1.056577346173265, 1.066722608119813, 1.020980703421571)
    9 = @(1.02, 1.044206523589279, 1.052563177271093, 1.051869958938809, 1.062018466855965, 1.043447689491414, 1.0499549077716, 1.055672653924425, 1.054981621361359, 1.065098441416134, 1.020455761612548)
    10 = @(1.02, 1.04263654304556, 1.051531519388439, 1.050524050613492, 1.060656821375043, 1.043085936983683, 1.048923021859836, 1.054922528283693, 1.053918546625454, 1.064016576921716, 1.02010505131308)
    11 = @(1.02, 1.041957327980737, 1.051085278042012, 1.049942235149113, 1.060068234084552, 1.042928121710911, 1.048476075635133, 1.054597344306137, 1.053458410559253, 1.063548347445032, 1.01995301651804)
    12 = @(1.02, 1.04170512740409, 1.050919596633123, 1.049726270077934, 1.059849760210192, 1.042869326160401, 1.04831004065582, 1.054476501144992, 1.053287523781272, 1.063374460999481, 1.019896518077131)
    13 = @(1.02, 1.041759221197277, 1.050955132554942, 1.049772588644542, 1.059896616624196, 1.042881945959494, 1.048345656582343, 1.054502424903175, 1.053324178283372, 1.063411758639675, 1.019908638353406)
    14 = @(1.02, 1.04193647916083, 1.051071581278171, 1.049924380401802, 1.06005017182407, 1.042923265238665, 1.04846235152059, 1.054587356496544, 1.053444284423746, 1.063533973220727, 1.019948346866719)
    15 = @(1.02, 1.042045705601014, 1.051143338827653, 1.050017923838221, 1.060144802649433, 1.042948700111331, 1.04853424859357, 1.054639678347324, 1.053518289549108, 1.063609278303788, 1.019972809163777)
    16 = @(1.02, 1.042681632927597, 1.051561145047563, 1.05056268435726, 1.06069590545223, 1.043096385966597, 1.048952681470445, 1.054944101859783, 1.053949088254706, 1.064047656597982, 1.020115137704594)
    17 = @(1.02, 1.043080693465847, 1.051823351551407, 1.050904659177439, 1.061041870112208, 1.043188711396717, 1.049215118296109, 1.055134959045161, 1.054219366313421, 1.064322700681374, 1.020204370040405)
    18 = @(1.02, 1.043313516283934, 1.051976337745153, 1.051104220987059, 1.061243763213576, 1.04324244994859, 1.049368180394382, 1.055246246604069, 1.054377032465158, 1.064483150880978, 1.020256400852623)
    19 = @(1.02, 1.04339291259688, 1.052028509772075, 1.051172282263551, 1.061312620087995, 1.043260754172948, 1.049420368411444, 1.055284186615883, 1.05443079549634, 1.064537863891752, 1.020274139138459)
    20 = @(1.02, 1.043037872089863, 1.051795214550866, 1.050867958817573, 1.061004741256434, 1.043178817468651, 1.049186962633664, 1.055114485613883, 1.054190366215901, 1.064293188808529, 1.020194798000857)
    21 = @(1.02, 1.041884278616497, 1.051037288019234, 1.049879677443084, 1.060004949393119, 1.0429111025937, 1.048427988295588, 1.054562347779705, 1.053408915372411, 1.063497983102985, 1.019936654412906)
    22 = @(1.02, 1.041159490375138, 1.050561169650577, 1.049259157328878, 1.059377230186269, 1.042741761968507, 1.047950679974359, 1.054214877013993, 1.052917750050221, 1.062998207739744, 1.019774199326596)
    23 = @(1.02, 1.041543664720577, 1.050813528721627, 1.049588025689072, 1.059709911260813, 1.042831628947209, 1.04820372043545, 1.054399107843482, 1.053178110264794, 1.063263128622614, 1.019860333946051)
    24 = @(1.02, 1.043057221048566, 1.051807928309598, 1.050884541841076, 1.061021517886198, 1.043183288459007, 1.049199685005351, 1.055123736788176, 1.054203470057562, 1.064306523886771, 1.020199123245746)
    25 = @(1.02, 1.044815875072951, 1.052963656505506, 1.052392747622107, 1.062547394964673, 1.043586923184081, 1.050354939330049, 1.055672653924425, 1.054981621361359, 1.065098441416134, 1.020455761612548)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $ws.Cells.Item($row, $colIndexes[$i]).Value = $vals[$i]
    }
}

Write-Host "Done updating rows 2-25"